$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ch_OldPassToInvalidPasswordTest")
$ws2 = $wb.Worksheets.Item("ch_OldPassToValidPasswordTest")

# --- Sheet "ch_OldPassToInvalidPasswordTest" (sheet1.xml) ---
# E3 value changes from 12 to 123456789
$ws1.Range("E3").Value = 123456789

# --- Sheet "ch_OldPassToValidPasswordTest" (sheet2.xml) ---
# Column E ("old password") is removed entirely; column F ("new password") shifts left into E.
$ws2.Columns("E").Delete()

# Update the (now shifted) values in column E for rows 3 and 4 to the new test data.
$ws2.Range("E3").Value = "112233445566"
$ws2.Range("E4").Value = "a"

# Update the selection on this (now inactive) sheet.
[void]$ws2.Range("E4").Select()

# Finally, make sheet1 the active / selected sheet+cell (must be last: selecting a
# range activates its sheet, so this also flips the workbook's active tab back).
$ws1.Activate()
[void]$ws1.Range("E3").Select()
